# Update the Price (D) / Volume(1h) (E) columns of the cryptos list with the
# latest refreshed values from the scheduled GitHub Actions data pull.
# D46/D49 use a leading apostrophe to force literal-text interpretation
# (preserving the trailing zero / significant digits) the same way a user
# typing '0.0360 into a cell would, since those strings would otherwise be
# auto-parsed as numbers and lose the trailing zero.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.864.13"
$ws.Range("E2").Value = "  +1.94%  "

$ws.Range("D3").Value = "3.033.95"
$ws.Range("E3").Value = "  +1.05%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "593.86"
$ws.Range("E5").Value = "  -1.03%  "

$ws.Range("D6").Value = "153.58"
$ws.Range("E6").Value = "  +5.77%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "3.030.54"
$ws.Range("E8").Value = "  +0.95%  "

$ws.Range("E9").Value = "  -1.13%  "

$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +11.17%  "

$ws.Range("E11").Value = "  +2.43%  "

$ws.Range("E12").Value = "  +0.98%  "

$ws.Range("E13").Value = "  +1.65%  "

$ws.Range("D14").Value = "35.62"
$ws.Range("E14").Value = "  +3.09%  "

$ws.Range("E15").Value = "  +1.97%  "

$ws.Range("D16").Value = "3.536.24"
$ws.Range("E16").Value = "  +1.06%  "

$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "62.817.02"
$ws.Range("E18").Value = "  +1.88%  "

$ws.Range("D19").Value = "3.035.20"

$ws.Range("D20").Value = "453.13"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").Value = "14.26"
$ws.Range("E21").Value = "  +1.53%  "

$ws.Range("E22").Value = "  +1.19%  "

$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("D24").Value = "83.14"
$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("D25").Value = "11.26"
$ws.Range("E25").Value = "  +2.63%  "

$ws.Range("E26").Value = "  +0.89%  "

$ws.Range("D27").Value = "12.27"
$ws.Range("E27").Value = "  +1.92%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "7.49"
$ws.Range("E29").Value = "  +2.83%  "

$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +7.91%  "

$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("D33").Value = "27.55"
$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("E34").Value = "  +0.48%  "

$ws.Range("D35").Value = "0.0₃0865"
$ws.Range("E35").Value = "  +3.26%  "

$ws.Range("E36").Value = "  +1.99%  "

$ws.Range("E37").Value = "  +2.44%  "

$ws.Range("D38").Value = "3.21"
$ws.Range("E38").Value = "  +10.26%  "

$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("E40").Value = "  +4.12%  "

$ws.Range("D41").Value = "50.38"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("D42").Value = "9.08"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("D43").Value = "0.304"
$ws.Range("E43").Value = "  +10.98%  "

$ws.Range("D44").Value = "42.22"
$ws.Range("E44").Value = "  +6.19%  "

$ws.Range("D45").Value = "393.67"
$ws.Range("E45").Value = "  -2.56%  "

$ws.Range("D46").Value = "'0.0360"

$ws.Range("D47").Value = "2.723.13"
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("D48").Value = "132.22"
$ws.Range("E48").Value = "  +0.25%  "

$ws.Range("D49").Value = "'2.30"
$ws.Range("E49").Value = "  +6.13%  "

$ws.Range("D51").Value = "24.41"
$ws.Range("E51").Value = "  +3.07%  "
